$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New task text shared by the two new entries
$newTask = "Sprint Planning Meeting 3"
$newDate = Get-Date -Year 2017 -Month 1 -Day 11 -Hour 0 -Minute 0 -Second 0
$newFrom = 0.59027777777777779   # 14:10:00
$newTo   = 0.65972222222222221   # 15:50:00

# Row 23: fill in the previously-empty "Sammer Manuel" block (columns W:AA)
$ws.Range("W23").Value = $newTask
$ws.Range("X23").Value = $newDate
$ws.Range("Y23").Value = $newFrom
$ws.Range("Z23").Value = $newTo
$ws.Range("AA23").Value = 1

# Row 28: fill in the previously-empty "Lamprecht Daniel" block (columns C:G)
# and "Ruhdorfer Alexander" block (columns M:Q)
$ws.Range("C28").Value = $newTask
$ws.Range("D28").Value = $newDate
$ws.Range("E28").Value = $newFrom
$ws.Range("F28").Value = $newTo
$ws.Range("G28").Value = 1

$ws.Range("M28").Value = $newTask
$ws.Range("N28").Value = $newDate
$ws.Range("O28").Value = $newFrom
$ws.Range("P28").Value = $newTo
$ws.Range("Q28").Value = 1

# Update the active selection to match the saved view state
$ws.Range("E36").Select()
